$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (D) and volume-change (E) values per row.
$updates = @(
    @{ Row = 2; D = '30.467.41'; E = '  +0.14%  ' }
    @{ Row = 3; D = '1.912.31'; E = '  -0.23%  ' }
    @{ Row = 4; D = '''0.9988'; E = '  -0.20%  ' }
    @{ Row = 5; D = '''244.72'; E = '  +1.46%  ' }
    @{ Row = 6; D = '''0.9988'; E = '  -0.18%  ' }
    @{ Row = 7; D = '''0.4814'; E = '  +2.30%  ' }
    @{ Row = 8; D = '''0.2887'; E = '  +1.31%  ' }
    @{ Row = 9; D = '''0.06716'; E = '  -1.29%  ' }
    @{ Row = 10; D = '''110.43'; E = '  +3.36%  ' }
    @{ Row = 11; D = '''19.03'; E = '  +4.33%  ' }
    @{ Row = 12; D = '1.911.89'; E = '  -0.16%  ' }
    @{ Row = 13; D = '''0.07546'; E = '  -1.64%  ' }
    @{ Row = 14; D = '''5.259'; E = '  +1.09%  ' }
    @{ Row = 15; D = '''0.6705'; E = '  +2.21%  ' }
    @{ Row = 16; D = '''288.85'; E = '  +0.01%  ' }
    @{ Row = 17; D = '30.471.73'; E = '  +0.10%  ' }
    @{ Row = 18; D = '''0.9990'; E = '  -0.16%  ' }
    @{ Row = 19; D = '''0.000007580'; E = '  -0.56%  ' }
    @{ Row = 20; D = '''12.83'; E = '  -0.80%  ' }
    @{ Row = 21; D = '2.163.02'; E = '  +0.57%  ' }
    @{ Row = 22; D = '''5.481'; E = '  +4.98%  ' }
    @{ Row = 23; D = '''0.9987'; E = '  -0.28%  ' }
    @{ Row = 24; D = '''6.405'; E = '  +3.43%  ' }
    @{ Row = 25; D = '''9.444'; E = '  +1.75%  ' }
    @{ Row = 26; D = '''164.23'; E = '  -2.30%  ' }
    @{ Row = 27; D = '''20.34'; E = '  -5.44%  ' }
    @{ Row = 28; D = '''2.114'; E = '  +2.29%  ' }
    @{ Row = 29; D = '''0.1053'; E = '  -1.59%  ' }
    @{ Row = 30; D = '''1.403'; E = '  +2.37%  ' }
    @{ Row = 31; D = '''4.163'; E = '  +0.15%  ' }
    @{ Row = 32; D = '''4.036'; E = '  +1.84%  ' }
    @{ Row = 33; D = '''0.04976'; E = '  -1.59%  ' }
    @{ Row = 34; D = '''0.7286'; E = '  -1.89%  ' }
    @{ Row = 35; D = '''1.132'; E = '  -1.68%  ' }
    @{ Row = 36; D = '''0.9992'; E = '  -0.08%  ' }
    @{ Row = 37; D = '''2.720'; E = '  -0.90%  ' }
    @{ Row = 38; D = '''0.02033'; E = '  -2.71%  ' }
    @{ Row = 39; D = '''2.667'; E = '  -0.70%  ' }
    @{ Row = 40; D = '''110.51'; E = '  +1.57%  ' }
    @{ Row = 41; D = '''2.015'; E = '  -1.81%  ' }
    @{ Row = 42; D = '''0.4438'; E = '  +4.83%  ' }
    @{ Row = 43; D = '''0.8651'; E = '  -0.35%  ' }
    @{ Row = 44; D = '''5.783'; E = '  -1.28%  ' }
    @{ Row = 45; D = '''0.9987'; E = '  -0.13%  ' }
    @{ Row = 46; D = '''68.14'; E = '  +0.99%  ' }
    @{ Row = 47; D = '''7.322'; E = '  +2.22%  ' }
    @{ Row = 48; D = '''49.10'; E = '  -1.85%  ' }
    @{ Row = 49; D = '''9.244'; E = '  +0.51%  ' }
    @{ Row = 50; D = '''0.1239'; E = '  +2.36%  ' }
    @{ Row = 51; D = '''34.79'; E = '  -0.10%  ' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
